$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ContactData")

# Update the mainImageURI value in C2 (LS LOGO.png -> DP.png)
$ws.Range("C2").Value = "https://rmoosa2014.github.io/Resume/DP.png"

# Match the saved selection state (activeCell E11)
$ws.Range("E11").Select()
